# "Finished the equipment section."
#  - Fixed a few issues.
#  - Adventures can now give XP bonuses.
#
# Adds 5 new affix rows (17-21) to the Affixes sheet, each referencing a
# newly-introduced name/description shared string pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 - Dark Hopes (suffix)
$ws.Range("A17").Value = "Dark Hopes"
$ws.Range("B17").Value = "The darkest of dreams and hopes imbue this item for all shadows will hide deep in the hearts of men and women."
$ws.Range("C17").Value = 0.15
$ws.Range("D17").Value = 0.15
$ws.Range("E17").Value = 0.15
$ws.Range("F17").Value = 0.13
$ws.Range("G17").Value = 0.13
$ws.Range("H17").Value = 0.13
$ws.Range("I17").Value = 0.13
$ws.Range("J17").Value = 0.13
$ws.Range("K17").Value = 23
$ws.Range("L17").Value = 12
$ws.Range("M17").Value = 40
$ws.Range("Q17").Value = 5000
$ws.Range("R17").Value = "suffix"

# Row 18 - Festering Doubt (prefix)
$ws.Range("A18").Value = "Festering Doubt"
$ws.Range("B18").Value = "Let it fester, let the doubt take over and let it make you run from battle."
$ws.Range("C18").Value = 0.05
$ws.Range("D18").Value = 0.05
$ws.Range("E18").Value = 0.05
$ws.Range("G18").Value = 0.3
$ws.Range("H18").Value = 0.3
$ws.Range("K18").Value = 24
$ws.Range("L18").Value = 15
$ws.Range("M18").Value = 40
$ws.Range("Q18").Value = 5000
$ws.Range("R18").Value = "prefix"

# Row 19 - Treasures Winds (prefix, Looting skill)
$ws.Range("A19").Value = "Treasures Winds"
$ws.Range("B19").Value = "Follow the winds to the treasure."
$ws.Range("K19").Value = 28
$ws.Range("L19").Value = 16
$ws.Range("M19").Value = 40
$ws.Range("N19").Value = "Looting"
$ws.Range("O19").Value = 0.15
$ws.Range("P19").Value = 0.05
$ws.Range("Q19").Value = 5000
$ws.Range("R19").Value = "prefix"

# Row 20 - Swfit Beat (suffix, Dodge skill)
$ws.Range("A20").Value = "Swfit Beat"
$ws.Range("B20").Value = "The music is upbeat and the drums pound with the sounds of war. Move quickly my child."
$ws.Range("K20").Value = 28
$ws.Range("L20").Value = 18
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = "Dodge"
$ws.Range("O20").Value = 0.15
$ws.Range("P20").Value = 0.05
$ws.Range("Q20").Value = 5000
$ws.Range("R20").Value = "suffix"

# Row 21 - Take Aim (prefix, Accuracy skill)
$ws.Range("A21").Value = "Take Aim"
$ws.Range("B21").Value = "Take aim at the enemy. Truly, they are vile. Never miss your mark, srtike em dead the first time. Right through the eyes."
$ws.Range("K21").Value = 29
$ws.Range("L21").Value = 20
$ws.Range("M21").Value = 40
$ws.Range("N21").Value = "Accuracy"
$ws.Range("O21").Value = 0.15
$ws.Range("P21").Value = 0.05
$ws.Range("Q21").Value = 5000
$ws.Range("R21").Value = "prefix"
